# Allow for text formatting
# - Replace the worked example (rows 2-9 with A/B/C/.../H process steps)
#   with a single generic "ExampleProcess" / "DEFER()" row.
# - Shrink the used range down to A1:B2.
# - Give column A an explicit width so text wraps/formats nicely next to
#   the wide column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 with the new example content.
$ws.Range("A2").Value2 = "ExampleProcess"
$ws.Range("B2").Value2 = "DEFER()"

# Drop the old rows 3-9 (A/B/C/D/E/F/G/H steps) entirely, shrinking the
# sheet's dimension down to A1:B2 and collapsing the shared-string table.
$ws.Range("A3:B9").EntireRow.Delete() | Out-Null

# Give column A a real width (previously only column B had a custom width).
$ws.Columns.Item(1).ColumnWidth = 15.15

# Move the selection, matching the new cursor position recorded in the file.
$ws.Range("B12").Select() | Out-Null
